$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 2645
$ws.Range("F9").Value = 530
$ws.Range("F11").Value = 1534
$ws.Range("F13").Value = 23
$ws.Range("F14").Value = 632
$ws.Range("F15").Value = 1519
$ws.Range("F16").Value = 1379
$ws.Range("F17").Value = 21
$ws.Range("F18").Value = 4
$ws.Range("F19").Value = 544
$ws.Range("F20").Value = 3905
$ws.Range("F21").Value = 3905
$ws.Range("F22").Value = 655
$ws.Range("F23").Value = 3304
$ws.Range("F24").Value = 761
$ws.Range("F25").Value = 24
$ws.Range("F26").Value = 2198
$ws.Range("F27").Value = 32
$ws.Range("F28").Value = 307
$ws.Range("F30").Value = 27
$ws.Range("F31").Value = 1165
$ws.Range("F32").Value = 764
$ws.Range("F34").Value = 1046
$ws.Range("F35").Value = 1044
$ws.Range("F36").Value = 78

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 93
$ws.Range("F7").Value = 16
$ws.Range("F11").Value = 10
$ws.Range("F18").Value = 252

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 242
$ws.Range("F3").Value = 384
$ws.Range("F4").Value = 532
$ws.Range("F5").Value = 137

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 242
$ws.Range("F6").Value = 384
$ws.Range("F9").Value = 532
$ws.Range("F10").Value = 2645
$ws.Range("F11").Value = 2645
$ws.Range("F17").Value = 93
$ws.Range("F18").Value = 16
$ws.Range("F19").Value = 530
$ws.Range("F23").Value = 1534
$ws.Range("F24").Value = 10
$ws.Range("F25").Value = 23
$ws.Range("F26").Value = 1519
$ws.Range("F28").Value = 1379
$ws.Range("F29").Value = 21
$ws.Range("F30").Value = 544
$ws.Range("F32").Value = 3905
$ws.Range("F33").Value = 3905
$ws.Range("F34").Value = 655
$ws.Range("F35").Value = 3304
$ws.Range("F36").Value = 761
$ws.Range("F37").Value = 2198
$ws.Range("F38").Value = 32
$ws.Range("F39").Value = 307
$ws.Range("F41").Value = 27
$ws.Range("F42").Value = 1165
$ws.Range("F44").Value = 252
$ws.Range("F47").Value = 764
$ws.Range("F49").Value = 1046
$ws.Range("F50").Value = 1044
$ws.Range("F51").Value = 78
